$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.242.41'
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").Value = '3.155.76'
$ws.Range("E3").Value = '  -0.62%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = "'570.20"
$ws.Range("E5").Value = '  +0.11%  '
$ws.Range("D6").Value = '163.29'
$ws.Range("E6").Value = '  -3.42%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").Value = '0.576'
$ws.Range("E8").Value = '  -5.27%  '
$ws.Range("E9").Value = '  -3.27%  '
$ws.Range("D10").Value = '6.59'
$ws.Range("E10").Value = '  -1.40%  '
$ws.Range("D11").Value = '0.381'
$ws.Range("E11").Value = '  -0.88%  '
$ws.Range("D12").Value = '3.704.15'
$ws.Range("E12").Value = '  -0.52%  '
$ws.Range("E13").Value = '  -0.79%  '
$ws.Range("D14").Value = '64.361.19'
$ws.Range("E14").Value = '  +0.17%  '
$ws.Range("D15").Value = "'25.10"
$ws.Range("E15").Value = '  -1.04%  '
$ws.Range("D16").Value = '3.150.81'
$ws.Range("E16").Value = '  -0.45%  '
$ws.Range("D17").Value = '0.0000154'
$ws.Range("E17").Value = '  -2.54%  '
$ws.Range("D18").Value = '404.26'
$ws.Range("E18").Value = '  -3.16%  '
$ws.Range("D19").Value = '12.65'
$ws.Range("E19").Value = '  -1.17%  '
$ws.Range("D20").Value = '5.22'
$ws.Range("E20").Value = '  -2.41%  '
$ws.Range("D21").Value = "'7.10"
$ws.Range("E21").Value = '  +0.30%  '
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("D23").Value = "'68.50"
$ws.Range("E23").Value = '  -2.18%  '
$ws.Range("D24").Value = '0.483'
$ws.Range("E24").Value = '  -1.55%  '
$ws.Range("E25").Value = '  -4.48%  '
$ws.Range("D26").Value = '0.0000101'
$ws.Range("E26").Value = '  -4.34%  '
$ws.Range("E27").Value = '  +0.32%  '
$ws.Range("D28").Value = '0.997'
$ws.Range("E28").Value = '  -1.45%  '
$ws.Range("E29").Value = '  -1.03%  '
$ws.Range("D30").Value = '21.14'
$ws.Range("E30").Value = '  -3.01%  '
$ws.Range("D31").Value = '6.26'
$ws.Range("E31").Value = '  -1.25%  '
$ws.Range("D32").Value = '4.81'
$ws.Range("E32").Value = '  -4.09%  '
$ws.Range("D33").Value = '156.72'
$ws.Range("E33").Value = '  +0.93%  '
$ws.Range("E34").Value = '  -2.07%  '
$ws.Range("D35").Value = '1.33'
$ws.Range("E35").Value = '  -2.98%  '
$ws.Range("D36").Value = '2.667.31'
$ws.Range("E36").Value = '  -1.18%  '
$ws.Range("D37").Value = '1.67'
$ws.Range("E37").Value = '  -1.77%  '
$ws.Range("D38").Value = '23.73'
$ws.Range("E38").Value = '  -3.44%  '
$ws.Range("D39").Value = '4.07'
$ws.Range("E39").Value = '  -2.29%  '
$ws.Range("D40").Value = '0.694'
$ws.Range("E40").Value = '  -2.27%  '
$ws.Range("E41").Value = '  -1.23%  '
$ws.Range("E42").Value = '  -3.97%  '
$ws.Range("D43").Value = '0.0255'
$ws.Range("E43").Value = '  -2.58%  '
$ws.Range("D44").Value = '288.01'
$ws.Range("E44").Value = '  -2.80%  '
$ws.Range("D45").Value = '21.22'
$ws.Range("E45").Value = '  -3.08%  '
$ws.Range("E46").Value = '  +0.08%  '
$ws.Range("D47").Value = '0.0979'
$ws.Range("E47").Value = '  -1.44%  '
$ws.Range("D48").Value = '10.51'
$ws.Range("E48").Value = '  +0.64%  '
$ws.Range("D49").Value = '1.89'
$ws.Range("E49").Value = '  -8.02%  '
$ws.Range("D50").Value = '5.68'
$ws.Range("E50").Value = '  -1.79%  '
$ws.Range("D51").Value = '0.871'
$ws.Range("E51").Value = '  -6.91%  '
